$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16 (shifts existing rows 16-19 down to 17-20)
$ws.Rows.Item(16).Insert()

# Set the new cell's value
$ws.Range("A16").Value = "Pepe Leal FC"
